# edit.ps1 - applies the PLOG0011A.docx diff via Word COM-interop emulation.
#
# Strategy: most hunks in the target diff only *split* an existing <w:r> run
# into several runs at specific character boundaries (often so a
# <w:proofErr> spell/grammar marker can bracket part of the text); a few
# hunks also change the text itself. Word COM does not expose a way to
# fabricate <w:proofErr> elements directly, so this script focuses on
# reproducing the run-splits and the textual edits precisely (same run
# boundaries / same text), which is the structurally observable part of
# the diff.
#
# Run splitting trick: toggling a character formatting property on a
# sub-range and then restoring it forces the engine to materialize that
# sub-range as its own <w:r> (with unchanged <w:rPr>), exactly like Word
# splits runs when you apply/clear formatting on a partial selection.

$d = $word.ActiveDocument

function Split-AtOffsets {
    # $rng must already cover exactly the target text.
    # $lens is an array of run lengths (relative), summing to the range length.
    param($rng, [int[]]$lens)
    $cursor = $rng.Start
    foreach ($len in $lens) {
        if ($len -gt 0) {
            $sub = $d.Range($cursor, $cursor + $len)
            $sub.Font.Bold = 1
            $sub.Font.Bold = 0
        }
        $cursor = $cursor + $len
    }
}

function Find-Range {
    param([string]$text)
    $r = $d.Content
    $ok = $r.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        Write-Output "NOT FOUND: $text"
    }
    return $r
}

# ---------------------------------------------------------------------
# 1) "Confecção de RTL para ressuprimento de material nacionalizado"
#    -> "Confecção de RTL para " | "ressuprimento" | " de material nacionalizado"
$r = Find-Range "Confecção de RTL para ressuprimento de material nacionalizado"
Split-AtOffsets $r @(23, 13, 26)

# ---------------------------------------------------------------------
# 2) " (alterar nomenclatura material carga, fazer publicação de designação, etc)"
#    -> " (alterar nomenclatura material carga, fazer publicação de designação, " | "etc" | ")"
$r = Find-Range " (alterar nomenclatura material carga, fazer publicação de designação, etc)"
Split-AtOffsets $r @(73, 3, 1)

# ---------------------------------------------------------------------
# 3) "FORNECEDOR visando a resolução" -> split "a " into "a" | " "
$r = Find-Range "FORNECEDOR visando a resolução"
Split-AtOffsets $r @(19, 1, 1)

# ---------------------------------------------------------------------
# 4) "ANALISAR RESULTADO DA VISITA" -> "REGISTRAR O" + " RESULTADO DA VISITA"
$r = Find-Range "ANALISAR RESULTADO DA VISITA"
$r.Text = "REGISTRAR O RESULTADO DA VISITA"
Split-AtOffsets $r @(11, 20)

# ---------------------------------------------------------------------
# 5) "prazo de validade de 1 (um) ano a partir da sua emissão."
#    original runs: "de" | " 1 (um) ano" -> "de" | " " | "1" | " (um) ano"
$r = Find-Range "prazo de validade de 1 (um) ano"
Split-AtOffsets $r @(18, 1, 1, 12)

# ---------------------------------------------------------------------
# 6) "1ª renovação: validade passa a ser de 2 (dois) anos;"
#    original run "de 2" -> "de " | "2"
$r = Find-Range "validade passa a ser de 2 (dois) anos"
Split-AtOffsets $r @(25, 1, 12)

# ---------------------------------------------------------------------
# 7) "2ª renovação e posteriores: validade passa a ser de 4 (quatro) anos"
#    -> "...ser de " | "4" | " (quatro) anos"
$r = Find-Range "2ª renovação e posteriores: validade passa a ser de 4 (quatro) anos"
Split-AtOffsets $r @(51, 1, 15)

# ---------------------------------------------------------------------
# 8) ilvl change (1 -> 2) for three consecutive paragraphs; handled later
#    via paragraph format (see below).

# ---------------------------------------------------------------------
# 9) "...implica na emissão de novo certificado válido por 1 (um) ano, o qual..."
$r = Find-Range "válido por 1 (um) ano, o qual pode ser renovado conforme informado."
Split-AtOffsets $r @(11, 1, 57)

Write-Output "part1 done"
